$d = $word.ActiveDocument

$d.Content.Find.Execute("Nur Faizah Mas Mohd Khalik", $false, $false, $false, $false, $false,
                         $true, 1, $false, "Nur Faizah", 2)
